$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (GET /api/articles/:article_id/comments): fill in the "Normal" (C)
# and "Queries" (D) checklist columns as "done", matching the formatting
# already used elsewhere in the sheet (copy format from C5 which uses the
# same "done" style).
$ws.Range("C5").Copy()
$ws.Range("C12:D12").PasteSpecial(-4122)
$ws.Range("C12").Value = "done"
$ws.Range("D12").Value = "done"

# Rows 16 & 17 (/api/users, /api/users/:username): the "Normal" column was
# marked plain "done", update it to call out that it's done but not using
# the framework defaults - reuse the format already used for similar
# caveats (copy format from C10 which uses that style).
$ws.Range("C10").Copy()
$ws.Range("C16:C17").PasteSpecial(-4122)
$ws.Range("C16").Value = "done (not defaults)"
$ws.Range("C17").Value = "done (not defaults)"

# Update the active selection to reflect where work left off.
$ws.Range("D12").Select()
